# Insert a new weekly price record as row 361 on the "Zanahoria" sheet.
# Inserting the row shifts every existing row from 361..410 down to 362..411,
# which reproduces the rest of the diff (all the D/J/K/L/M/O/P "changes"
# further down are really just the pre-existing rows sliding down by one).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 361..410 down to 362..411, leaving row 361 free for the new record.
$ws.Rows.Item(361).Insert()

# Populate the newly inserted row 361 with the new weekly record.
$ws.Range("A361").Value = 4
$ws.Range("B361").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C361").Value = "Los Lagos"
$ws.Range("D361").Value = 44776
$ws.Range("E361").Value = 10
$ws.Range("F361").Value = 100114013
$ws.Range("G361").Value = "Zanahoria"
$ws.Range("H361").Value = "Sin especificar"
$ws.Range("I361").Value = "Primera"
$ws.Range("J361").Value = 150
$ws.Range("K361").Value = 10000
$ws.Range("L361").Value = 10000
$ws.Range("M361").Value = 10000
$ws.Range("N361").Value = "`$/saco 20 kilos"
$ws.Range("O361").Value = "Provincia de Llanquihue"
$ws.Range("P361").Value = 500
$ws.Range("Q361").Value = 20
$ws.Range("R361").Value = "Hortaliza"
